# Updates the cryptos list (coin prices + 1h volume change %) per the
# "Updated cryptos list" GitHub Actions commit.
#
# Source cells are plain text (coin price strings like "67.097.20" or
# "1.00" are *not* numbers -- commas/locale formatting and the mixed
# "thousand-separator" style mean they must stay text). Excel's COM
# `.Value` setter auto-coerces plain numeric-looking strings to actual
# numbers, which would both change the cell type and drop formatting
# (e.g. "1.00" -> 1). To avoid that we prefix such values with a
# leading apostrophe (the standard "force text" marker) and restore
# the cell's original Style afterwards, since the apostrophe marker
# otherwise nudges Excel to allocate a new (quote-prefixed) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
    $cell.Style = $origStyle
}

Set-CellText $ws 2 4 "67.097.20"
Set-CellText $ws 2 5 "  -1.16%  "
Set-CellText $ws 3 4 "2.474.51"
Set-CellText $ws 3 5 "  -2.31%  "
Set-CellText $ws 4 4 "1.00"
Set-CellText $ws 4 5 "  -0.09%  "
Set-CellText $ws 5 4 "583.03"
Set-CellText $ws 5 5 "  -1.32%  "
Set-CellText $ws 6 4 "169.19"
Set-CellText $ws 6 5 "  -2.60%  "
Set-CellText $ws 7 5 "  +0.01%  "
Set-CellText $ws 8 5 "  -1.94%  "
Set-CellText $ws 9 4 "2.474.13"
Set-CellText $ws 9 5 "  -2.33%  "
Set-CellText $ws 10 5 "  -2.51%  "
Set-CellText $ws 11 5 "  -0.83%  "
Set-CellText $ws 12 4 "4.93"
Set-CellText $ws 12 5 "  -2.12%  "
Set-CellText $ws 13 5 "  -3.82%  "
Set-CellText $ws 14 4 "25.53"
Set-CellText $ws 14 5 "  -3.35%  "
Set-CellText $ws 16 4 "67.018.91"
Set-CellText $ws 16 5 "  -1.29%  "
Set-CellText $ws 17 5 "  -4.27%  "
Set-CellText $ws 18 4 "2.485.69"
Set-CellText $ws 18 5 "  -1.68%  "
Set-CellText $ws 19 4 "11.07"
Set-CellText $ws 19 5 "  -5.96%  "
Set-CellText $ws 20 5 "  -5.95%  "
Set-CellText $ws 21 4 "351.49"
Set-CellText $ws 21 5 "  -4.92%  "
Set-CellText $ws 22 4 "4.05"
Set-CellText $ws 22 5 "  -2.21%  "
Set-CellText $ws 23 4 "0.998"
Set-CellText $ws 23 5 "  -0.20%  "
Set-CellText $ws 24 4 "68.90"
Set-CellText $ws 24 5 "  -4.18%  "
Set-CellText $ws 25 4 "4.26"
Set-CellText $ws 25 5 "  -6.75%  "
Set-CellText $ws 26 4 "1.83"
Set-CellText $ws 26 5 "  -4.62%  "
Set-CellText $ws 27 4 "9.20"
Set-CellText $ws 27 5 "  -7.51%  "
Set-CellText $ws 28 5 "  -58.95%  "
Set-CellText $ws 29 4 "2.596.18"
Set-CellText $ws 29 5 "  -2.65%  "
Set-CellText $ws 30 5 "  -6.61%  "
Set-CellText $ws 31 4 "510.00"
Set-CellText $ws 31 5 "  -5.33%  "
Set-CellText $ws 32 5 "  -7.20%  "
Set-CellText $ws 33 5 "  -6.44%  "
Set-CellText $ws 34 5 "  -5.29%  "
Set-CellText $ws 35 4 "0.999"
Set-CellText $ws 35 5 "  -0.09%  "
Set-CellText $ws 36 4 "159.49"
Set-CellText $ws 36 5 "  -0.19%  "
Set-CellText $ws 37 5 "  -10.07%  "
Set-CellText $ws 38 4 "18.66"
Set-CellText $ws 38 5 "  +0.29%  "
Set-CellText $ws 39 4 "18.34"
Set-CellText $ws 39 5 "  -4.40%  "
Set-CellText $ws 40 5 "  -7.62%  "
Set-CellText $ws 41 5 "  -4.81%  "
Set-CellText $ws 42 5 "  -0.01%  "
Set-CellText $ws 43 5 "  -6.26%  "
Set-CellText $ws 44 4 "0.328"
Set-CellText $ws 44 5 "  -6.43%  "
Set-CellText $ws 45 4 "2.37"
Set-CellText $ws 45 5 "  -6.29%  "
Set-CellText $ws 46 4 "38.96"
Set-CellText $ws 46 5 "  -1.05%  "
Set-CellText $ws 47 4 "141.06"
Set-CellText $ws 47 5 "  -4.63%  "
Set-CellText $ws 48 2 "Filecoin"
Set-CellText $ws 48 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws 48 4 "3.45"
Set-CellText $ws 48 5 "  -7.00%  "
Set-CellText $ws 49 2 "ARBITRUM"
Set-CellText $ws 49 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText $ws 49 4 "0.514"
Set-CellText $ws 49 5 "  -6.56%  "
Set-CellText $ws 50 4 "0.0₆0254"
Set-CellText $ws 50 5 "  -11.04%  "
Set-CellText $ws 51 5 "  -7.21%  "

Write-Output "cryptos list updated"
